{"js": "// Append, after the last paragraph of the document body, three new\n// paragraphs:\n//   1. an empty paragraph (spacer)\n//   2. a paragraph containing a hyperlink to the LearnOpenGL anti-aliasing\n//      article\n//   3. a paragraph with a short text comment about the link\nconst body = context.document.body;\n\n// Anchor on the very last paragraph currently in the document (the ImGUI\n// paragraph) and insert everything after it, in order.\nconst lastParagraph = body.paragraphs.getLast();\n\n// 1) Blank spacer paragraph.\nconst blankParagraph = lastParagraph.insertParagraph(\"\", \"After\");\n\n// 2) Paragraph whose entire contents is a hyperlink.\nconst linkParagraph = blankParagraph.insertParagraph(\"\", \"After\");\nconst linkUrl = \"https://learnopengl.com/Advanced-OpenGL/Anti-Aliasing\";\nconst linkRange = linkParagraph.insertText(linkUrl, \"Replace\");\nlinkRange.hyperlink = linkUrl;\n\n// 3) Paragraph with the trailing descriptive text.\nconst textParagraph = linkParagraph.insertParagraph(\"\", \"After\");\ntextParagraph.insertText(\n  \"This Website covers Anti Aliasing, not really needed but would be a nice addition to adopt early while developing the base rendering systems.\",\n  \"Replace\"\n);\n\nawait context.sync();\n", "ps1": "# Append, after the last paragraph of the document body, three new\n# paragraphs:\n#   1. an empty paragraph (spacer)\n#   2. a paragraph containing a hyperlink to the LearnOpenGL anti-aliasing\n#      article\n#   3. a paragraph with a short text comment about the link\n$d = $word.ActiveDocument\n\n# 1) Blank spacer paragraph, added right after the current last paragraph\n#    (the \"ImGUI will be a useful library...\" paragraph).\n$tail = $d.Paragraphs.Last.Range\n$tail.Collapse(0)   # wdCollapseEnd\n$tail.InsertParagraphAfter()\n\n# 2) Paragraph whose entire contents is a hyperlink to the article.\n$linkRange = $d.Paragraphs.Last.Range\n$linkRange.Collapse(0)\n$linkRange.InsertParagraphAfter()\n\n$linkUrl = \"https://learnopengl.com/Advanced-OpenGL/Anti-Aliasing\"\n$hyperlinkRange = $d.Paragraphs.Last.Range\n$hyperlinkRange.Collapse(0)\n$hyperlinkRange.InsertAfter($linkUrl)\n$d.Hyperlinks.Add($hyperlinkRange, $linkUrl)\n\n# 3) Paragraph with the trailing descriptive text.\n$textRange = $d.Paragraphs.Last.Range\n$textRange.Collapse(0)\n$textRange.InsertParagraphAfter()\n\n$finalRange = $d.Paragraphs.Last.Range\n$finalRange.Collapse(0)\n$finalRange.InsertAfter(\"This Website covers Anti Aliasing, not really needed but would be a nice addition to adopt early while developing the base rendering systems.\")\n"}
